$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6456.4585
$ws.Range("I40").Value = 7227.769
$ws.Range("K40").Value = 7227.769
$ws.Range("M40").Value = -7052.769
$ws.Range("H100").Value = 1720.7693
$ws.Range("I100").Value = 1261.25
$ws.Range("J100").Value = 2456
$ws.Range("K100").Value = 1261.25
$ws.Range("L100").Value = 2456
$ws.Range("M100").Value = -720.25
$ws.Range("N100").Value = -3538
$ws.Range("H137").Value = 4588.914
$ws.Range("I137").Value = 5233.5186
$ws.Range("J137").Value = 2413.375
$ws.Range("K137").Value = 15700.5558
$ws.Range("L137").Value = 7240.125
$ws.Range("M137").Value = -13150.5558
$ws.Range("N137").Value = -12340.125
$ws.Range("H138").Value = 5330.2
$ws.Range("I138").Value = 2767.2856
$ws.Range("J138").Value = 6011.481
$ws.Range("K138").Value = 8301.856800000001
$ws.Range("L138").Value = 18034.443
$ws.Range("M138").Value = -3161.856800000001
$ws.Range("N138").Value = -28314.443
$ws.Range("H139").Value = 289994.16
$ws.Range("J139").Value = 289994.16
$ws.Range("L139").Value = 289994.16
$ws.Range("N139").Value = -300274.16

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 206.625
$ws.Range("I5").Value = 211.14285
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 211.14285
$ws.Range("L5").Value = 175
$ws.Range("M5").Value = -99.14285000000001
$ws.Range("N5").Value = -399
$ws.Range("H32").Value = 12028.897
$ws.Range("I32").Value = 11636.806
$ws.Range("K32").Value = 11636.806
$ws.Range("M32").Value = -11349.806
$ws.Range("H45").Value = 2211.6924
$ws.Range("I45").Value = 1034.4445
$ws.Range("K45").Value = 1034.4445
$ws.Range("M45").Value = -657.4445000000001
$ws.Range("H61").Value = 3346.6978
$ws.Range("I61").Value = 2074.8708
$ws.Range("K61").Value = 2074.8708
$ws.Range("M61").Value = -1862.8708
$ws.Range("H74").Value = 5565.7964
$ws.Range("I74").Value = 5611.24
$ws.Range("K74").Value = 5611.24
$ws.Range("M74").Value = -4737.24
$ws.Range("H77").Value = 5565.7964
$ws.Range("I77").Value = 5611.24
$ws.Range("K77").Value = 28056.2
$ws.Range("M77").Value = -23688.2
$ws.Range("H110").Value = 15454.467
$ws.Range("I110").Value = 16179.5
$ws.Range("K110").Value = 16179.5
$ws.Range("M110").Value = -14134.5
$ws.Range("H136").Value = 3346.6978
$ws.Range("I136").Value = 2074.8708
$ws.Range("K136").Value = 6224.6124
$ws.Range("M136").Value = -3674.6124

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 206.625
$ws.Range("I4").Value = 211.14285
$ws.Range("J4").Value = 175
$ws.Range("K4").Value = 211.14285
$ws.Range("L4").Value = 175
$ws.Range("M4").Value = -96.14285000000001
$ws.Range("N4").Value = -405
$ws.Range("H80").Value = 802.7778
$ws.Range("I80").Value = 837.75
$ws.Range("J80").Value = 774.8
$ws.Range("K80").Value = 837.75
$ws.Range("L80").Value = 774.8
$ws.Range("M80").Value = 160.25
$ws.Range("N80").Value = -2770.8
$ws.Range("H82").Value = 90932536
$ws.Range("I82").Value = 100017590
$ws.Range("K82").Value = 100017590
$ws.Range("M82").Value = -100017207
$ws.Range("H83").Value = 802.7778
$ws.Range("I83").Value = 837.75
$ws.Range("J83").Value = 774.8
$ws.Range("K83").Value = 4188.75
$ws.Range("L83").Value = 3874
$ws.Range("M83").Value = 803.25
$ws.Range("N83").Value = -13858
$ws.Range("H85").Value = 90932536
$ws.Range("I85").Value = 100017590
$ws.Range("K85").Value = 100017590
$ws.Range("M85").Value = -100016264
$ws.Range("H94").Value = 1274.8462
$ws.Range("J94").Value = 1483.1666
$ws.Range("L94").Value = 1483.1666
$ws.Range("N94").Value = -2385.1666
$ws.Range("H134").Value = 2535.2297
$ws.Range("J134").Value = 2499.3635
$ws.Range("L134").Value = 7498.0905
$ws.Range("N134").Value = -12568.0905
$ws.Range("H138").Value = 59999.145
$ws.Range("J138").Value = 59999.145
$ws.Range("L138").Value = 59999.145
$ws.Range("N138").Value = -70279.14499999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9999.333000000001
$ws.Range("I25").Value = 10000
$ws.Range("J25").Value = 9999
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 9999
$ws.Range("M25").Value = -9826
$ws.Range("N25").Value = -10347
$ws.Range("H31").Value = 38921716
$ws.Range("I31").Value = 11114603
$ws.Range("J31").Value = 66728830
$ws.Range("K31").Value = 11114603
$ws.Range("L31").Value = 66728830
$ws.Range("M31").Value = -11114308
$ws.Range("N31").Value = -66729420
$ws.Range("H34").Value = 38921716
$ws.Range("I34").Value = 11114603
$ws.Range("J34").Value = 66728830
$ws.Range("K34").Value = 11114603
$ws.Range("L34").Value = 66728830
$ws.Range("M34").Value = -11114401
$ws.Range("N34").Value = -66729234
$ws.Range("H59").Value = 43750
$ws.Range("J59").Value = 43750
$ws.Range("L59").Value = 43750
$ws.Range("N59").Value = -46040
$ws.Range("H107").Value = 5865.8887
$ws.Range("I107").Value = 906
$ws.Range("K107").Value = 906
$ws.Range("M107").Value = 1014
$ws.Range("H141").Value = 130056.47
$ws.Range("J141").Value = 130056.47
$ws.Range("L141").Value = 130056.47
$ws.Range("N141").Value = -140416.47

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 91.44444
$ws.Range("I38").Value = 105.8
$ws.Range("K38").Value = 317.4
$ws.Range("M38").Value = 29.60000000000002
$ws.Range("H40").Value = 54.727272
$ws.Range("I40").Value = 59.625
$ws.Range("K40").Value = 238.5
$ws.Range("M40").Value = -169.5
$ws.Range("H55").Value = 6510
$ws.Range("J55").Value = 7887.5
$ws.Range("L55").Value = 23662.5
$ws.Range("N55").Value = -24016.5
$ws.Range("H113").Value = 2390.2942
$ws.Range("I113").Value = 1520.1111
$ws.Range("K113").Value = 4560.3333
$ws.Range("M113").Value = -2390.3333
$ws.Range("H122").Value = 1503.2727
$ws.Range("J122").Value = 1730.125
$ws.Range("L122").Value = 15571.125
$ws.Range("N122").Value = -20471.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 14002.333
$ws.Range("I11").Value = 11001.5
$ws.Range("K11").Value = 11001.5
$ws.Range("M11").Value = -10862.5
$ws.Range("H18").Value = 9999
$ws.Range("J18").Value = 9999
$ws.Range("L18").Value = 9999
$ws.Range("N18").Value = -10585
$ws.Range("H24").Value = 12498
$ws.Range("J24").Value = 12498
$ws.Range("L24").Value = 12498
$ws.Range("N24").Value = -12844
$ws.Range("H70").Value = 66672864
$ws.Range("I70").Value = 5498.8
$ws.Range("J70").Value = 100006540
$ws.Range("K70").Value = 5498.8
$ws.Range("L70").Value = 100006540
$ws.Range("M70").Value = -5228.8
$ws.Range("N70").Value = -100007080
$ws.Range("H73").Value = 66672864
$ws.Range("I73").Value = 5498.8
$ws.Range("J73").Value = 100006540
$ws.Range("K73").Value = 5498.8
$ws.Range("L73").Value = 100006540
$ws.Range("M73").Value = -4562.8
$ws.Range("N73").Value = -100008412
$ws.Range("H107").Value = 66598.336
$ws.Range("I107").Value = 92015
$ws.Range("K107").Value = 92015
$ws.Range("M107").Value = -90095
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0
$ws.Range("H122").Value = 2195.1538
$ws.Range("I122").Value = 2363.4285
$ws.Range("K122").Value = 7090.2855
$ws.Range("M122").Value = -4640.2855
$ws.Range("H132").Value = 37496.855
$ws.Range("I132").Value = 51547.35
$ws.Range("J132").Value = 2370.625
$ws.Range("K132").Value = 154642.05
$ws.Range("L132").Value = 7111.875
$ws.Range("M132").Value = -152112.05
$ws.Range("N132").Value = -12171.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("N3").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0
$ws.Range("H21").Value = 6211.3076
$ws.Range("I21").Value = 2375.25
$ws.Range("J21").Value = 7916.222
$ws.Range("K21").Value = 2375.25
$ws.Range("L21").Value = 7916.222
$ws.Range("M21").Value = -2201.25
$ws.Range("N21").Value = -8264.222
$ws.Range("H23").Value = 5673.778
$ws.Range("I23").Value = 3883.25
$ws.Range("K23").Value = 3883.25
$ws.Range("M23").Value = -3653.25
$ws.Range("H132").Value = 9271.406999999999
$ws.Range("I132").Value = 3482.6667
$ws.Range("J132").Value = 9995
$ws.Range("K132").Value = 10448.0001
$ws.Range("L132").Value = 29985
$ws.Range("M132").Value = -7918.000100000001
$ws.Range("N132").Value = -35045

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("H62").Value = 22738990
$ws.Range("J62").Value = 27789320
$ws.Range("L62").Value = 27789320
$ws.Range("N62").Value = -27790568
$ws.Range("H65").Value = 22738990
$ws.Range("J65").Value = 27789320
$ws.Range("L65").Value = 138946600
$ws.Range("N65").Value = -138952840
$ws.Range("H126").Value = 3460.4348
$ws.Range("I126").Value = 2189.2632
$ws.Range("J126").Value = 9498.5
$ws.Range("K126").Value = 6567.7896
$ws.Range("L126").Value = 28495.5
$ws.Range("M126").Value = -4097.7896
$ws.Range("N126").Value = -33435.5
